$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1. Update header row (row 1) with new snake_case column names
$ws.Range("A1").Value = "mx_state"
$ws.Range("B1").Value = "mx_municipality"
$ws.Range("C1").Value = "n_matriculas"
$ws.Range("D1").Value = "pct_matriculas"

# 2. Title-case the Spanish connector words (de/del/la/las/los/el/y) in
#    state/municipality names throughout column A and B
$renames = @(
    @('B8', 'Pabellón De Arteaga'),
    @('B9', 'Rincón De Romos'),
    @('B14', 'Playas De Rosarito'),
    @('B30', 'Amatenango De La Frontera'),
    @('B31', 'Amatenango Del Valle'),
    @('B37', 'Comitán De Domínguez'),
    @('B53', 'San Cristóbal De Las Casas'),
    @('B79', 'Coyame Del Sotol'),
    @('B88', 'Guadalupe Y Calvo'),
    @('B91', 'Hidalgo Del Parral'),
    @('B111', 'San Francisco De Borja'),
    @('B112', 'San Francisco Del Oro'),
    @('B119', 'Valle De Zaragoza'),
    @('B130', 'San Juan De Sabinas'),
    @('B141', 'Villa De Álvarez'),
    @('A143', 'Ciudad De México'),
    @('B147', 'Cuajimalpa De Morelos'),
    @('B162', 'Coneto De Comonfort'),
    @('B172', 'Nombre De Dios'),
    @('B182', 'San Juan Del Río'),
    @('B183', 'San Luis Del Cordero'),
    @('A190', 'Estado De México'),
    @('B190', 'Acambay De Ruíz Castañeda'),
    @('B192', 'Almoloya De Alquisiras'),
    @('B198', 'Atizapán De Zaragoza'),
    @('B204', 'Chapa De Mota'),
    @('B209', 'Coacalco De Berriozábal'),
    @('B213', 'Ecatepec De Morelos'),
    @('B217', 'Ixtapan De La Sal'),
    @('B228', 'Naucalpan De Juárez'),
    @('B236', 'San Felipe Del Progreso'),
    @('B238', 'San Simón De Guerrero'),
    @('B240', 'Soyaniquilpan De Juárez'),
    @('B249', 'Tenango Del Valle'),
    @('B259', 'Tlalnepantla De Baz'),
    @('B264', 'Valle De Chalco Solidaridad'),
    @('B265', 'Villa De Allende'),
    @('B276', 'San Miguel De Allende'),
    @('B277', 'Apaseo El Alto'),
    @('B278', 'Apaseo El Grande'),
    @('B285', 'Dolores Hidalgo Cuna De La Independencia Nacional'),
    @('B289', 'Jaral Del Progreso'),
    @('B297', 'Purísima Del Rincón'),
    @('B301', 'San Diego De La Unión'),
    @('B303', 'San Francisco Del Rincón'),
    @('B305', 'San Luis De La Paz'),
    @('B306', 'Santa Cruz De Juventino Rosas'),
    @('B308', 'Silao De La Victoria'),
    @('B312', 'Valle De Santiago'),
    @('B317', 'Acapulco De Juárez'),
    @('B319', 'Ajuchitlán Del Progreso'),
    @('B323', 'Atenango Del Río'),
    @('B324', 'Atoyac De Álvarez'),
    @('B325', 'Ayutla De Los Libres'),
    @('B328', 'Buenavista De Cuéllar'),
    @('B329', 'Chilapa De Álvarez'),
    @('B330', 'Chilpancingo De Los Bravo'),
    @('B331', 'Coahuayutla De José María Izazaga'),
    @('B336', 'Coyuca De Benítez'),
    @('B337', 'Coyuca De Catalán'),
    @('B341', 'Cuetzala Del Progreso'),
    @('B342', 'Cutzamala De Pinzón'),
    @('B348', 'Huitzuco De Los Figueroa'),
    @('B349', 'Iguala De La Independencia'),
    @('B351', 'Ixcateopan De Cuauhtémoc'),
    @('B352', 'Zihuatanejo De Azueta'),
    @('B354', 'La Unión De Isidoro Montes De Oca'),
    @('B356', 'Mártir De Cuilapan'),
    @('B369', 'Taxco De Alarcón'),
    @('B371', 'Técpan De Galeana'),
    @('B373', 'Tepecoacuilco De Trujano'),
    @('B375', 'Tixtla De Guerrero'),
    @('B378', 'Tlalixtaquilla De Maldonado'),
    @('B379', 'Tlapa De Comonfort'),
    @('B388', 'Agua Blanca De Iturbide'),
    @('B393', 'Atotonilco El Grande'),
    @('B397', 'Cuautepec De Hinojosa'),
    @('B399', 'Huasca De Ocampo'),
    @('B403', 'Huejutla De Reyes'),
    @('B406', 'Jacala De Ledezma'),
    @('B410', 'Mineral De La Reforma'),
    @('B411', 'Mineral Del Monte'),
    @('B412', 'Mixquiahuala De Juárez'),
    @('B413', 'Molango De Escamilla'),
    @('B414', 'Pachuca De Soto'),
    @('B423', 'Tepeji Del Río De Ocampo'),
    @('B425', 'Tezontepec De Aldama'),
    @('B433', 'Tula De Allende'),
    @('B434', 'Tulancingo De Bravo'),
    @('B435', 'Villa De Tezontepec'),
    @('B439', 'Zacualtipán De Ángeles'),
    @('B444', 'Acatlán De Juárez'),
    @('B445', 'Ahualulco De Mercado'),
    @('B450', 'Atemajac De Brizuela'),
    @('B453', 'Atotonilco El Alto'),
    @('B455', 'Autlán De Navarro'),
    @('B461', 'Cañadas De Obregón'),
    @('B466', 'Concepción De Buenos Aires'),
    @('B467', 'Cuautitlán De García Barragán'),
    @('B478', 'Huejuquilla El Alto'),
    @('B479', 'Ixtlahuacán De Los Membrillos'),
    @('B480', 'Ixtlahuacán Del Río'),
    @('B483', 'Jilotlán De Los Dolores'),
    @('B489', 'La Manzanilla De La Paz'),
    @('B490', 'Lagos De Moreno'),
    @('B497', 'Ojuelos De Jalisco'),
    @('B502', 'San Cristóbal De La Barranca'),
    @('B503', 'San Diego De Alejandría'),
    @('B505', 'San Juan De Los Lagos'),
    @('B506', 'San Juanito De Escobedo'),
    @('B509', 'San Miguel El Alto'),
    @('B510', 'San Sebastián Del Oeste'),
    @('B511', 'Santa María De Los Ángeles'),
    @('B514', 'Talpa De Allende'),
    @('B515', 'Tamazula De Gordiano'),
    @('B521', 'Teocuitatlán De Corona'),
    @('B522', 'Tepatitlán De Morelos'),
    @('B525', 'Tizapán El Alto'),
    @('B526', 'Tlajomulco De Zúñiga'),
    @('B537', 'Unión De San Antonio'),
    @('B538', 'Unión De Tula'),
    @('B539', 'Valle De Guadalupe'),
    @('B540', 'Valle De Juárez'),
    @('B545', 'Yahualica De González Gallo'),
    @('B546', 'Zacoalco De Torres'),
    @('B549', 'Zapotitlán De Vadillo'),
    @('B550', 'Zapotlán Del Rey'),
    @('B551', 'Zapotlán El Grande'),
    @('B573', 'Coalcomán De Vázquez Pallares'),
    @('B635', 'Tiquicheo De Nicolás Romero'),
    @('B657', 'Coatlán Del Río'),
    @('B665', 'Jonacatepec De Leandro Valle'),
    @('B669', 'Puente De Ixtla'),
    @('B675', 'Tlaltizapán De Zapata'),
    @('B680', 'Zacualpan De Amilpas'),
    @('B684', 'Amatlán De Cañas'),
    @('B685', 'Bahía De Banderas'),
    @('B689', 'Ixtlán Del Río'),
    @('B694', 'Santa María Del Oro'),
    @('B707', 'Mier Y Noriega'),
    @('B710', 'San Nicolás De Los Garza'),
    @('B714', 'Acatlán De Pérez Figueroa'),
    @('B715', 'Constancia Del Rosario'),
    @('B717', 'Guadalupe De Ramírez'),
    @('B718', 'Heroica Ciudad De Ejutla De Crespo'),
    @('B719', 'Heroica Ciudad De Huajuapan De León'),
    @('B720', 'Heroica Ciudad De Tlaxiaco'),
    @('B721', 'Ixtlán De Juárez'),
    @('B722', 'Heroica Ciudad De Juchitán De Zaragoza'),
    @('B726', 'Mártires De Tacubaya'),
    @('B728', 'Miahuatlán De Porfirio Díaz'),
    @('B729', 'Oaxaca De Juárez'),
    @('B730', 'Ocotlán De Morelos'),
    @('B731', 'Pinotepa De Don Luis'),
    @('B733', 'Putla Villa De Guerrero'),
    @('B734', 'Reforma De Pineda'),
    @('B742', 'San Baltazar Yatzachi El Bajo'),
    @('B761', 'San Juan Bautista Lo De Soto'),
    @('B784', 'San Pablo Villa De Mitla'),
    @('B797', 'Santa Cruz Tacache De Mina'),
    @('B800', 'Santa Inés Del Monte'),
    @('B807', 'Santa María Jalapa Del Marqués'),
    @('B833', 'Tamazulápam Del Espíritu Santo'),
    @('B834', 'Teotitlán De Flores Magón'),
    @('B835', 'Teotitlán Del Valle'),
    @('B836', 'Heroica Villa Tezoatlán De Segura Y Luna, Cuna De La Independencia De Oaxaca'),
    @('B837', 'Tlacolula De Matamoros'),
    @('B838', 'Totontepec Villa De Morelos'),
    @('B839', 'Villa De Etla'),
    @('B840', 'Villa De Tututepec'),
    @('B841', 'Villa De Zaachila'),
    @('B843', 'Villa Sola De Vega'),
    @('B844', 'Zimatlán De Álvarez'),
    @('B870', 'Cuayuca De Andrade'),
    @('B876', 'Huehuetlán El Grande'),
    @('B878', 'Ixcamilpa De Guerrero'),
    @('B881', 'Izúcar De Matamoros'),
    @('B887', 'Los Reyes De Juárez'),
    @('B892', 'Palmar De Bravo'),
    @('B907', 'San Salvador El Seco'),
    @('B912', 'Tecali De Herrera'),
    @('B919', 'Tepexi De Rodríguez'),
    @('B924', 'Tlacotepec De Benito Juárez'),
    @('B942', 'Amealco De Bonfil'),
    @('B944', 'Cadereyta De Montes'),
    @('B947', 'Jalpan De Serra'),
    @('B950', 'Pinal De Amoles'),
    @('B952', 'San Juan Del Río'),
    @('B965', 'Ciudad Del Maíz'),
    @('B972', 'Mexquitic De Carmona'),
    @('B977', 'San Ciro De Acosta'),
    @('B980', 'Santa María Del Río'),
    @('B984', 'Villa De Arriaga'),
    @('B985', 'Villa De Guadalupe'),
    @('B986', 'Villa De La Paz'),
    @('B987', 'Villa De Ramos'),
    @('B988', 'Villa De Reyes'),
    @('B1027', 'Nacozari De García'),
    @('B1068', 'Ixtacuixtla De Mariano Matamoros'),
    @('B1069', 'Mazatecochco De José María Morelos'),
    @('B1072', 'Papalotla De Xicohténcatl'),
    @('B1074', 'Sanctórum De Lázaro Cárdenas'),
    @('B1076', 'Tepetitla De Lardizábal'),
    @('B1079', 'Tetla De La Solidaridad'),
    @('B1088', 'Ziltlaltépec De Trinidad Sánchez Santos'),
    @('B1097', 'Alto Lucero De Gutiérrez Barrios'),
    @('B1100', 'Amatlán De Los Reyes'),
    @('B1107', 'Boca Del Río'),
    @('B1109', 'Castillo De Teayo'),
    @('B1111', 'Cazones De Herrera'),
    @('B1120', 'Cosamaloapan De Carpio'),
    @('B1121', 'Cosautlán De Carvajal'),
    @('B1132', 'Ignacio De La Llave'),
    @('B1134', 'Ixhuacán De Los Reyes'),
    @('B1135', 'Ixhuatlán De Madero'),
    @('B1143', 'Juchique De Ferrer'),
    @('B1148', 'Martínez De La Torre'),
    @('B1149', 'Medellín De Bravo'),
    @('B1162', 'Paso De Ovejas'),
    @('B1163', 'Paso Del Macho'),
    @('B1166', 'Poza Rica De Hidalgo'),
    @('B1173', 'Sayula De Alemán'),
    @('B1175', 'Soledad De Doblado'),
    @('B1194', 'Vega De Alatorre'),
    @('B1201', 'Zozocolco De Hidalgo'),
    @('B1227', 'Cañitas De Felipe Pescador'),
    @('B1236', 'Jiménez Del Teul'),
    @('B1241', 'Mezquital Del Oro'),
    @('B1246', 'Nochistlán De Mejía'),
    @('B1247', 'Noria De Ángeles'),
    @('B1258', 'Teúl De González Ortega'),
    @('B1259', 'Tlaltenango De Sánchez Román'),
    @('B1261', 'Trinidad García De La Cadena'),
    @('B1263', 'Villa De Cos')
)
foreach ($pair in $renames) {
    $ws.Range($pair[0]).Value = $pair[1]
}

# 3. Fix the floating point representation of the 11/11486 percentage
#    (recomputed value differs from the original by 1 ULP)
$floatRows = @(8,272,275,352,358,396,456,460,463,514,568,571,578,689,837,904,912,999,1044,1070,1104,1150,1166,1258)
foreach ($r in $floatRows) {
    $ws.Range("D" + $r).Value = 0.0009576876197109524
}

# 4. Remove the trailing footer/metadata rows (1271-1276), shrinking the
#    used range down to A1:D1270
$ws.Range("A1271:A1276").EntireRow.Delete() | Out-Null

Write-Host ("Final dimension: " + $ws.UsedRange.Address())